$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "A little worse"
$ws.Range("C3").Value = "Somewhat worse"
$ws.Range("C4").Value = "Somewhat worse"
$ws.Range("C5").Value = "Somewhat worse"
$ws.Range("C6").Value = "A lot worse"
$ws.Range("C7").Value = "A little worse"
$ws.Range("C8").Value = "Not worse"
$ws.Range("C9").Value = "A lot worse"
$ws.Range("C10").Value = "Somewhat worse"
$ws.Range("C11").Value = "Somewhat worse"
$ws.Range("C12").Value = "Somewhat worse"
$ws.Range("C13").Value = "Somewhat worse"
$ws.Range("C14").Value = "Somewhat worse"
$ws.Range("C15").Value = "A little worse"

[void]$ws.Range("C16").Select()
